$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.502.32'
$ws.Range('D2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.305.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.51%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.19%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E9').Value = '  -3.98%  '

$ws.Range('E10').Value = '  -2.28%  '

$ws.Range('E11').Value = '  -5.03%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.877.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.48%  '

$ws.Range('E13').Value = '  -0.92%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.81%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '66.521.10'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.54%  '

$ws.Range('E16').Value = '  -2.98%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.309.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.75%  '

$ws.Range('E18').Value = '  -0.91%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.85%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.12%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.38%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '

$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.517'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000118'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.41%  '

$ws.Range('E26').Value = '  +0.31%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.93%  '

$ws.Range('E28').Value = '  -0.56%  '

$ws.Range('E29').Value = '  -2.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.78%  '

$ws.Range('E31').Value = '  -5.63%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.01%  '

$ws.Range('E33').Value = '  -3.93%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.15%  '

$ws.Range('E35').Value = '  -1.52%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.94%  '

$ws.Range('E37').Value = '  -4.67%  '

$ws.Range('E38').Value = '  -0.61%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.802.03'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.49%  '

$ws.Range('E40').Value = '  -3.28%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.43'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.75%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0673'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.47%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.58%  '

$ws.Range('E45').Value = '  -4.84%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '320.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.09%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0272'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.41%  '

$ws.Range('E49').Value = '  -3.06%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.25%  '

$ws.Range('E51').Value = '  -1.46%  '
